$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add note on 'In Touch' - new row 19 entry: word "festival", pronunciation, part of speech "n."
$ws.Range("A19").Value = "festival"
$ws.Range("B19").Value = "/'festIvl/"
$ws.Range("C19").Value = "n."

# Update the active selection to B20
$ws.Range("B20").Select()
